$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.001.49'
$ws.Range('E2').Value = '  -3.81%  '
$ws.Range('D3').Value = '1.642.07'
$ws.Range('E3').Value = '  -5.78%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9968'
$ws.Range('E4').Value = '  -0.38%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.49'
$ws.Range('E5').Value = '  -5.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9989'
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4798'
$ws.Range('E7').Value = '  -4.87%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2600'
$ws.Range('E8').Value = '  -5.24%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06110'
$ws.Range('E9').Value = '  -1.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07077'
$ws.Range('E10').Value = '  -2.35%  '
$ws.Range('D11').Value = '1.640.29'
$ws.Range('E11').Value = '  -5.89%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.59'
$ws.Range('E12').Value = '  -3.69%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6012'
$ws.Range('E13').Value = '  -7.72%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.387'
$ws.Range('E14').Value = '  -6.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '73.79'
$ws.Range('E15').Value = '  -4.80%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9990'
$ws.Range('E16').Value = '  -0.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9981'
$ws.Range('E17').Value = '  -0.28%  '
$ws.Range('D18').Value = '25.002.20'
$ws.Range('E18').Value = '  -3.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006617'
$ws.Range('E19').Value = '  -3.88%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.28'
$ws.Range('E20').Value = '  -5.15%  '
$ws.Range('D21').Value = '1.846.73'
$ws.Range('E21').Value = '  -6.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.396'
$ws.Range('E22').Value = '  -1.39%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.596'
$ws.Range('E23').Value = '  -1.50%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.255'
$ws.Range('E24').Value = '  -2.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '133.31'
$ws.Range('E25').Value = '  -1.82%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '14.94'
$ws.Range('E26').Value = '  -2.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.386'
$ws.Range('E27').Value = '  -7.99%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '104.35'
$ws.Range('E28').Value = '  -1.16%  '
$ws.Range('E29').Value = '  -7.29%  '
$ws.Range('E30').Value = '  -0.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07705'
$ws.Range('E31').Value = '  -5.86%  '
$ws.Range('E32').Value = '  -2.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.9977'
$ws.Range('E33').Value = '  -0.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04291'
$ws.Range('E34').Value = '  -7.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.573'
$ws.Range('E35').Value = '  -3.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9288'
$ws.Range('E36').Value = '  -6.68%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5894'
$ws.Range('E37').Value = '  -3.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.565'
$ws.Range('E38').Value = '  -7.70%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01522'
$ws.Range('E39').Value = '  -5.97%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9977'
$ws.Range('E40').Value = '  -0.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8307'
$ws.Range('E41').Value = '  +8.66%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '98.54'
$ws.Range('E42').Value = '  -2.13%  '
$ws.Range('E43').Value = '  -8.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3700'
$ws.Range('E44').Value = '  -5.36%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.688'
$ws.Range('E45').Value = '  -6.30%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1094'
$ws.Range('E46').Value = '  -5.91%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.095'
$ws.Range('E47').Value = '  -3.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05199'
$ws.Range('E48').Value = '  -2.17%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '29.20'
$ws.Range('E49').Value = '  -4.74%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.9984'
$ws.Range('E50').Value = '  -0.44%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9971'
$ws.Range('E51').Value = '  -0.37%  '
